$wb = $excel.ActiveWorkbook

# --- paper_summary: trim trailing space in "- Male / - Weight " note ---
$wsSummary = $wb.Worksheets.Item("paper_summary")
$wsSummary.Range("B5").Value = "- Male`n- Weight"

# --- control_data: trim trailing spaces off several header labels ---
$wsControl = $wb.Worksheets.Item("control_data")
$wsControl.Range("P1").Value = "BAL Macrophages"
$wsControl.Range("V1").Value = "BAL LDH"
$wsControl.Range("Y1").Value = "BAL Total Protein"

# --- dataset: fill in the newly-recorded "No. of Subjects (N)" column Q for rows 10-25 ---
$wsDataset = $wb.Worksheets.Item("dataset")
$wsDataset.Range("Q10:Q25").Value = 1

# Remove the stale reviewer comment left on F1 of the dataset sheet (also drops
# the associated legacy VML drawing reference).
$commentCell = $wsDataset.Range("F1")
if ($commentCell.Comment -ne $null) {
    $commentCell.Comment.Delete() | Out-Null
}

# --- sheet/tab selection state: dataset becomes the active sheet/tab instead
# of surface_area_calculation, with Q25 as the active selection ---
$wsDataset.Activate() | Out-Null
$wsDataset.Range("Q25").Select() | Out-Null
